$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4537.92
$ws.Range("I40").Value = 4406.727
$ws.Range("K40").Value = 4406.727
$ws.Range("M40").Value = -4231.727

$ws.Range("H113").Value = 3596.9167
$ws.Range("I113").Value = 3111.4
$ws.Range("J113").Value = 3943.7144
$ws.Range("K113").Value = 3111.4
$ws.Range("L113").Value = 3943.7144
$ws.Range("M113").Value = 142.5999999999999
$ws.Range("N113").Value = -10451.7144

$ws.Range("H116").Value = 3983.3235
$ws.Range("I116").Value = 3563.8
$ws.Range("K116").Value = 3563.8
$ws.Range("M116").Value = -121.8000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4094.6453
$ws.Range("I32").Value = 3598.1333
$ws.Range("J32").Value = 18990
$ws.Range("K32").Value = 3598.1333
$ws.Range("L32").Value = 18990
$ws.Range("M32").Value = -3311.1333
$ws.Range("N32").Value = -19564

$ws.Range("H39").Value = 4998
$ws.Range("I39").Value = 4998
$ws.Range("K39").Value = 4998
$ws.Range("M39").Value = -4478

$ws.Range("H74").Value = 1567
$ws.Range("I74").Value = 1362.2858
$ws.Range("K74").Value = 1362.2858
$ws.Range("M74").Value = -488.2858000000001

$ws.Range("H77").Value = 1567
$ws.Range("I77").Value = 1362.2858
$ws.Range("K77").Value = 6811.429
$ws.Range("M77").Value = -2443.429

$ws.Range("H122").Value = 11527.286
$ws.Range("I122").Value = 1944.1482
$ws.Range("K122").Value = 5832.444600000001
$ws.Range("M122").Value = -3382.444600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 640.65625
$ws.Range("I94").Value = 590.0714
$ws.Range("J94").Value = 994.75
$ws.Range("K94").Value = 590.0714
$ws.Range("L94").Value = 994.75
$ws.Range("M94").Value = -139.0714
$ws.Range("N94").Value = -1896.75

$ws.Range("H134").Value = 6429.6875
$ws.Range("I134").Value = 5067.769
$ws.Range("K134").Value = 15203.307
$ws.Range("M134").Value = -12668.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1966.0312
$ws.Range("I31").Value = 2014.2759
$ws.Range("J31").Value = 1926.0571
$ws.Range("K31").Value = 2014.2759
$ws.Range("L31").Value = 1926.0571
$ws.Range("M31").Value = -1719.2759
$ws.Range("N31").Value = -2516.0571

$ws.Range("H34").Value = 1966.0312
$ws.Range("I34").Value = 2014.2759
$ws.Range("J34").Value = 1926.0571
$ws.Range("K34").Value = 2014.2759
$ws.Range("L34").Value = 1926.0571
$ws.Range("M34").Value = -1812.2759
$ws.Range("N34").Value = -2330.0571

$ws.Range("H58").Value = 1569.5938
$ws.Range("I58").Value = 1563.7587
$ws.Range("J58").Value = 1626
$ws.Range("K58").Value = 1563.7587
$ws.Range("L58").Value = 1626
$ws.Range("M58").Value = -1360.7587
$ws.Range("N58").Value = -2032

$ws.Range("H93").Value = 11449
$ws.Range("I93").Value = 11449
$ws.Range("K93").Value = 11449
$ws.Range("M93").Value = -9577

$ws.Range("H136").Value = 1569.5938
$ws.Range("I136").Value = 1563.7587
$ws.Range("J136").Value = 1626
$ws.Range("K136").Value = 4691.2761
$ws.Range("L136").Value = 4878
$ws.Range("M136").Value = -2141.2761
$ws.Range("N136").Value = -9978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 290.64285
$ws.Range("I33").Value = 140.4
$ws.Range("J33").Value = 374.1111
$ws.Range("K33").Value = 842.4000000000001
$ws.Range("L33").Value = 2244.6666
$ws.Range("M33").Value = -559.4000000000001
$ws.Range("N33").Value = -2810.6666

$ws.Range("H97").Value = 641.25
$ws.Range("I97").Value = 369.66666
$ws.Range("J97").Value = 731.7778
$ws.Range("K97").Value = 1108.99998
$ws.Range("L97").Value = 2195.3334
$ws.Range("M97").Value = -612.9999800000001
$ws.Range("N97").Value = -3187.3334

$ws.Range("H122").Value = 1230.32
$ws.Range("J122").Value = 1378
$ws.Range("L122").Value = 12402
$ws.Range("N122").Value = -17302

$ws.Range("H128").Value = 198411.05
$ws.Range("I128").Value = 198411.05
$ws.Range("K128").Value = 595233.1499999999
$ws.Range("M128").Value = -590253.1499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 59999
$ws.Range("J88").Value = 59999
$ws.Range("L88").Value = 59999
$ws.Range("N88").Value = -60901

$ws.Range("H91").Value = 59999
$ws.Range("J91").Value = 59999
$ws.Range("L91").Value = 59999
$ws.Range("N91").Value = -63119

$ws.Range("H113").Value = 4140.0938
$ws.Range("I113").Value = 3962.1052
$ws.Range("J113").Value = 4400.231
$ws.Range("K113").Value = 3962.1052
$ws.Range("L113").Value = 4400.231
$ws.Range("M113").Value = -1792.1052
$ws.Range("N113").Value = -8740.231

$ws.Range("H122").Value = 2061.8064
$ws.Range("I122").Value = 1881.9524
$ws.Range("K122").Value = 5645.857199999999
$ws.Range("M122").Value = -3195.857199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6579592
$ws.Range("I16").Value = 7813114
$ws.Range("J16").Value = 810.6667
$ws.Range("K16").Value = 7813114
$ws.Range("L16").Value = 810.6667
$ws.Range("M16").Value = -7812944
$ws.Range("N16").Value = -1150.6667

$ws.Range("H39").Value = 18019.334
$ws.Range("I39").Value = 18019.334
$ws.Range("K39").Value = 18019.334
$ws.Range("M39").Value = -17559.334

$ws.Range("H42").Value = 7020
$ws.Range("I42").Value = 8768
$ws.Range("K42").Value = 8768
$ws.Range("M42").Value = -8205

$ws.Range("H45").Value = 43999.5
$ws.Range("I45").Value = 43999.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 43999.5
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -43592.5

$ws.Range("H46").Value = 2368.7144
$ws.Range("I46").Value = 906.375
$ws.Range("K46").Value = 906.375
$ws.Range("M46").Value = -718.375

$ws.Range("H47").Value = 30000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 30000
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("M47").Value = 30000
$ws.Range("N47").Value = -30980

$ws.Range("H48").Value = 33900
$ws.Range("I48").Value = 25850
$ws.Range("J48").Value = 50000
$ws.Range("K48").Value = 25850
$ws.Range("L48").Value = 50000
$ws.Range("M48").Value = -25189
$ws.Range("N48").Value = -51322

$ws.Range("H49").Value = 7020
$ws.Range("I49").Value = 8768
$ws.Range("K49").Value = 8768
$ws.Range("M49").Value = -8621

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0

$ws.Range("H52").Value = 30000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = 30000
$ws.Range("N52").Value = -30466

$ws.Range("H53").Value = 40000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 40000
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("M53").Value = 40000
$ws.Range("N53").Value = -41036

$ws.Range("H54").Value = 80000
$ws.Range("J54").Value = 80000
$ws.Range("L54").Value = 80000
$ws.Range("N54").Value = -81288

$ws.Range("H55").Value = 467.17648
$ws.Range("I55").Value = 322.1111
$ws.Range("J55").Value = 630.375
$ws.Range("K55").Value = 322.1111
$ws.Range("L55").Value = 630.375
$ws.Range("M55").Value = -149.1111
$ws.Range("N55").Value = -976.375

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H58").Value = 57165.5
$ws.Range("I58").Value = 35748.25
$ws.Range("K58").Value = 35748.25
$ws.Range("M58").Value = -35488.25

$ws.Range("H93").Value = 8335040.5
$ws.Range("I93").Value = 10527791
$ws.Range("J93").Value = 2585.8
$ws.Range("K93").Value = 10527791
$ws.Range("L93").Value = 2585.8
$ws.Range("M93").Value = -10526543
$ws.Range("N93").Value = -5081.8

$ws.Range("H99").Value = 49999.5
$ws.Range("I99").Value = 35000
$ws.Range("K99").Value = 35000
$ws.Range("M99").Value = -32005

$ws.Range("H100").Value = 250003920
$ws.Range("I100").Value = 250003920
$ws.Range("K100").Value = 250003920
$ws.Range("M100").Value = -250003379

$ws.Range("H132").Value = 13484.981
$ws.Range("I132").Value = 16713.342
$ws.Range("J132").Value = 2454.75
$ws.Range("K132").Value = 50140.026
$ws.Range("L132").Value = 7364.25
$ws.Range("M132").Value = -47610.026
$ws.Range("N132").Value = -12424.25

$ws.Range("H136").Value = 4288577.5
$ws.Range("I136").Value = 6002075
$ws.Range("J136").Value = 4833.0835
$ws.Range("K136").Value = 18006225
$ws.Range("L136").Value = 14499.2505
$ws.Range("M136").Value = -18003675
$ws.Range("N136").Value = -19599.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8267754.5
$ws.Range("I81").Value = 10103900
$ws.Range("K81").Value = 20207800
$ws.Range("M81").Value = -20206739

$ws.Range("H84").Value = 8267754.5
$ws.Range("I84").Value = 10103900
$ws.Range("K84").Value = 101039000
$ws.Range("M84").Value = -101033696

$ws.Range("H126").Value = 3684.919
$ws.Range("I126").Value = 3490.3333
$ws.Range("J126").Value = 4210.3
$ws.Range("K126").Value = 10470.9999
$ws.Range("L126").Value = 12630.9
$ws.Range("M126").Value = -8000.999899999999
$ws.Range("N126").Value = -17570.9

$ws.Range("H132").Value = 2161.7632
$ws.Range("I132").Value = 1892.1526
$ws.Range("K132").Value = 5676.4578
$ws.Range("M132").Value = -3146.4578
